# Update the "dSF" column (column F) values on the active worksheet.
# The source data was repulled, so the dSF (change in closing line vs. the
# game-closing line at game time) values differ from the dS0 values that
# were previously (incorrectly) copied into column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    5  = -1
    6  = 2
    7  = 2
    8  = -1
    10 = 5
    11 = 13
    13 = -2
    14 = -1
    15 = 1
    16 = 7
    17 = 1
    18 = 1
    19 = -1
    20 = -1
    21 = -1
    22 = 1
    23 = 6
    24 = -1
    25 = 8
    26 = -4
    27 = -2
    28 = -1
    29 = -2
    30 = -4
    33 = 6
    34 = -1
    35 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
